$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, $cellRef, $newValue)
    if ($null -eq $newValue) {
        $ws.Range($cellRef).Value = ""
    } else {
        $ws.Range($cellRef).Value = [double]$newValue
    }
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws "H57" 70840
Set-CellValue $ws "J57" 71680
Set-CellValue $ws "L57" 215040
Set-CellValue $ws "N57" -216038
Set-CellValue $ws "H98" 1517.875
Set-CellValue $ws "I98" 1517.875
Set-CellValue $ws "K98" 1517.875
Set-CellValue $ws "M98" -19.875
Set-CellValue $ws "H122" 1517.875
Set-CellValue $ws "I122" 1517.875
Set-CellValue $ws "K122" 4553.625
Set-CellValue $ws "M122" -2103.625
Set-CellValue $ws "H125" 1069160.9
Set-CellValue $ws "I125" 1832037.6
Set-CellValue $ws "K125" 16488338.4
Set-CellValue $ws "M125" -16485878.4
Set-CellValue $ws "H132" 1047.1555
Set-CellValue $ws "I132" 1047.1555
Set-CellValue $ws "K132" 3141.4665
Set-CellValue $ws "M132" -611.4665000000005
Set-CellValue $ws "H135" 968.1316
Set-CellValue $ws "I135" 913.2162
Set-CellValue $ws "K135" 8218.9458
Set-CellValue $ws "M135" -5683.9458
Set-CellValue $ws "H136" 184001
Set-CellValue $ws "J136" 184001
Set-CellValue $ws "L136" 184001
Set-CellValue $ws "N136" -194201
Set-CellValue $ws "H137" 5370.0415
Set-CellValue $ws "I137" 8021.4546
Set-CellValue $ws "K137" 24064.3638
Set-CellValue $ws "M137" -21514.3638
Set-CellValue $ws "H138" 9603.954
Set-CellValue $ws "H141" 1652.6
Set-CellValue $ws "I141" 1647.7241
Set-CellValue $ws "K141" 4943.1723
Set-CellValue $ws "M141" 236.8276999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws "H32" 22558.916
Set-CellValue $ws "I32" 12336.75
Set-CellValue $ws "J32" 43003.25
Set-CellValue $ws "K32" 12336.75
Set-CellValue $ws "L32" 43003.25
Set-CellValue $ws "M32" -12049.75
Set-CellValue $ws "N32" -43577.25
Set-CellValue $ws "H45" 5196.88
Set-CellValue $ws "I45" 4598.647
Set-CellValue $ws "K45" 4598.647
Set-CellValue $ws "M45" -4221.647
Set-CellValue $ws "H74" 3363.6
Set-CellValue $ws "I74" 3250.7273
Set-CellValue $ws "K74" 3250.7273
Set-CellValue $ws "M74" -2376.7273
Set-CellValue $ws "H77" 3363.6
Set-CellValue $ws "I77" 3250.7273
Set-CellValue $ws "K77" 16253.6365
Set-CellValue $ws "M77" -11885.6365

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws "H59" 0
Set-CellValue $ws "J59" 0
Set-CellValue $ws "L59" 0
Set-CellValue $ws "N59" $null
Set-CellValue $ws "H86" 1477.7778
Set-CellValue $ws "I86" 1462.75
Set-CellValue $ws "J86" 1489.8
Set-CellValue $ws "K86" 1462.75
Set-CellValue $ws "L86" 1489.8
Set-CellValue $ws "M86" -339.75
Set-CellValue $ws "N86" -3735.8
Set-CellValue $ws "H89" 1477.7778
Set-CellValue $ws "I89" 1462.75
Set-CellValue $ws "J89" 1489.8
Set-CellValue $ws "K89" 7313.75
Set-CellValue $ws "L89" 7449
Set-CellValue $ws "M89" -1697.75
Set-CellValue $ws "N89" -18681
Set-CellValue $ws "H105" 2159.2222
Set-CellValue $ws "I105" 2159.2222
Set-CellValue $ws "K105" 2159.2222
Set-CellValue $ws "M105" -412.2222000000002
Set-CellValue $ws "H134" 5302.7646
Set-CellValue $ws "I134" 3262.4167
Set-CellValue $ws "J134" 10199.6
Set-CellValue $ws "K134" 9787.250100000001
Set-CellValue $ws "L134" 30598.8
Set-CellValue $ws "M134" -7252.250100000001
Set-CellValue $ws "N134" -35668.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws "H52" 125518
Set-CellValue $ws "J52" 125518
Set-CellValue $ws "L52" 125518
Set-CellValue $ws "N52" -126106
Set-CellValue $ws "H99" 6139.7
Set-CellValue $ws "I99" 5466.3335
Set-CellValue $ws "K99" 5466.3335
Set-CellValue $ws "M99" -3968.3335
Set-CellValue $ws "H126" 6139.7
Set-CellValue $ws "I126" 5466.3335
Set-CellValue $ws "K126" 16399.0005
Set-CellValue $ws "M126" -13929.0005
Set-CellValue $ws "H132" 221591.83
Set-CellValue $ws "I132" 281869.4
Set-CellValue $ws "K132" 845608.2000000001
Set-CellValue $ws "M132" -843078.2000000001
Set-CellValue $ws "H141" 242041
Set-CellValue $ws "J141" 269547.84
Set-CellValue $ws "L141" 269547.84
Set-CellValue $ws "N141" -279907.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws "H4" 13376232
Set-CellValue $ws "I4" 10140750
Set-CellValue $ws "J4" 19598310
Set-CellValue $ws "K4" 30422250
Set-CellValue $ws "L4" 58794930
Set-CellValue $ws "M4" -30422138
Set-CellValue $ws "N4" -58795154
Set-CellValue $ws "H80" 15000
Set-CellValue $ws "J80" 15000
Set-CellValue $ws "L80" 45000
Set-CellValue $ws "N80" -46872
Set-CellValue $ws "H83" 15000
Set-CellValue $ws "J83" 15000
Set-CellValue $ws "L83" 135000
Set-CellValue $ws "N83" -144360
Set-CellValue $ws "H92" 1146.2
Set-CellValue $ws "I92" 211.33333
Set-CellValue $ws "J92" 2548.5
Set-CellValue $ws "K92" 633.99999
Set-CellValue $ws "L92" 7645.5
Set-CellValue $ws "M92" 614.00001
Set-CellValue $ws "N92" -10141.5
Set-CellValue $ws "H114" 5810.3335
Set-CellValue $ws "I114" 0
Set-CellValue $ws "J114" 5810.3335
Set-CellValue $ws "K114" 0
Set-CellValue $ws "L114" 17431.0005
Set-CellValue $ws "M114" $null
Set-CellValue $ws "N114" -23939.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws "H97" 1243.2858
Set-CellValue $ws "I97" 1173.5555
Set-CellValue $ws "J97" 1661.6666
Set-CellValue $ws "K97" 1173.5555
Set-CellValue $ws "L97" 1661.6666
Set-CellValue $ws "M97" -677.5554999999999
Set-CellValue $ws "N97" -2653.6666
Set-CellValue $ws "H102" 4663.25
Set-CellValue $ws "I102" 2151.125
Set-CellValue $ws "K102" 2151.125
Set-CellValue $ws "M102" -529.125
Set-CellValue $ws "H125" 100000
Set-CellValue $ws "J125" 100000
Set-CellValue $ws "L125" 100000
Set-CellValue $ws "N125" -104920

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws "H46" 79579.73
Set-CellValue $ws "J46" 79579.73
Set-CellValue $ws "L46" 79579.73
Set-CellValue $ws "N46" -80041.73
Set-CellValue $ws "H126" 6989.2
Set-CellValue $ws "I126" 4665.8335
Set-CellValue $ws "K126" 13997.5005
Set-CellValue $ws "M126" -11527.5005
Set-CellValue $ws "H132" 282197.7
Set-CellValue $ws "I132" 373671.66
Set-CellValue $ws "J132" 7775.778
Set-CellValue $ws "K132" 1121014.98
Set-CellValue $ws "L132" 23327.334
Set-CellValue $ws "M132" -1118484.98
Set-CellValue $ws "N132" -28387.334
Set-CellValue $ws "H134" 79579.73
Set-CellValue $ws "J134" 79579.73
Set-CellValue $ws "L134" 238739.19
Set-CellValue $ws "N134" -243809.19
